$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# 72×15=1080 -> 39×46=1794
$t.Cell(1, 1).Range.Text = "39×46=1794"

# 32×25=800 -> 76×36=2736
$t.Cell(1, 2).Range.Text = "76×36=2736"

# 18×27=486 -> 85×36=3060
$t.Cell(1, 3).Range.Text = "85×36=3060"

# 86×14=1204 -> 93×22=2046
$t.Cell(1, 4).Range.Text = "93×22=2046"

# 32×35=1120 -> 72×91=6552
$t.Cell(1, 5).Range.Text = "72×91=6552"

# 45×56=2520 -> 71×85=6035
$t.Cell(5, 1).Range.Text = "71×85=6035"

# 81×29=2349 -> 35×17=595
$t.Cell(5, 2).Range.Text = "35×17=595"

# 93×22=2046 -> 95×83=7885
$t.Cell(5, 3).Range.Text = "95×83=7885"

# 78×86=6708 -> 60×24=1440
$t.Cell(5, 4).Range.Text = "60×24=1440"

# 63×97=6111 -> 38×20=760
$t.Cell(5, 5).Range.Text = "38×20=760"

# 73×22=1606 -> 59×58=3422
$t.Cell(10, 1).Range.Text = "59×58=3422"

# 41×66=2706 -> 93×91=8463
$t.Cell(10, 2).Range.Text = "93×91=8463"

# 44×42=1848 -> 86×80=6880
$t.Cell(10, 3).Range.Text = "86×80=6880"

# 83×26=2158 -> 94×20=1880
$t.Cell(10, 4).Range.Text = "94×20=1880"

# 97×85=8245 -> 59×44=2596
$t.Cell(10, 5).Range.Text = "59×44=2596"

# 62×56=3472 -> 62×59=3658
$t.Cell(15, 1).Range.Text = "62×59=3658"

# 91×14=1274 -> 24×80=1920
$t.Cell(15, 2).Range.Text = "24×80=1920"

# 47×54=2538 -> 70×11=770
$t.Cell(15, 3).Range.Text = "70×11=770"

# 75×33=2475 -> 39×70=2730
$t.Cell(15, 4).Range.Text = "39×70=2730"

# 92×38=3496 -> 82×50=4100
$t.Cell(15, 5).Range.Text = "82×50=4100"

# 29×87=2523 -> 25×15=375
$t.Cell(20, 1).Range.Text = "25×15=375"

# 63×31=1953 -> 36×37=1332
$t.Cell(20, 2).Range.Text = "36×37=1332"

# 69×57=3933 -> 91×45=4095
$t.Cell(20, 3).Range.Text = "91×45=4095"

# 83×31=2573 -> 13×26=338
$t.Cell(20, 4).Range.Text = "13×26=338"

# 97×82=7954 -> 73×27=1971
$t.Cell(20, 5).Range.Text = "73×27=1971"

Write-Output "Replacements applied successfully"